$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Numeric data refresh: headers (row 1) and the recomputed statistics table,
#    plus a brand-new column H appended to the right of the existing data.
$numericValues = [ordered]@{
    "B1" = 1000
    "C1" = 1250
    "D1" = 1500
    "E1" = 1750
    "F1" = 2000
    "G1" = 2250
    "H1" = 2500
    "B2" = 0.019222546161321671
    "C2" = 0.019358358795285902
    "D2" = 0.019314641744548288
    "E2" = 0.019281045751633988
    "F2" = 0.019228221265319641
    "G2" = 0.019287399716847571
    "H2" = 0.019203177043300025
    "B3" = 0.077035957240038866
    "C3" = 0.077346137058053252
    "D3" = 0.07707820954254796
    "E3" = 0.077178649237472766
    "F3" = 0.076995693938390189
    "G3" = 0.076592732420953269
    "H3" = 0.077056110684089171
    "B4" = 0.30767735665694851
    "C4" = 0.30886075949367087
    "D4" = 0.30857517625840303
    "E4" = 0.30904139433551198
    "F4" = 0.30771778734680355
    "G4" = 0.30740915526191603
    "H4" = 0.30758390981296441
    "B5" = 0.69582118561710393
    "C5" = 0.69489305979921434
    "D5" = 0.69421216592884072
    "E5" = 0.69389978213507619
    "F5" = 0.6947664789665452
    "G5" = 0.69372345445965078
    "H5" = 0.69344094286446323
    "B6" = 1
    "C6" = 1
    "D6" = 1
    "E6" = 1
    "F6" = 1
    "G6" = 1
    "H6" = 1
    "B7" = 1.2361516034985423
    "C7" = 1.2295940637276299
    "D7" = 1.2338088211182159
    "E7" = 1.2358387799564272
    "F7" = 1.225405763497847
    "G7" = 1.2302973100519112
    "H7" = 1.2343069433768896
    "B8" = 1.9241982507288629
    "C8" = 1.9262330859886514
    "D8" = 1.9331037875061485
    "E8" = 1.9379084967320261
    "F8" = 1.9261344816164292
    "G8" = 1.9296838131193961
    "H8" = 1.9267230335639254
    "B9" = 2.7648202137998052
    "C9" = 2.7708424268878216
    "D9" = 2.7643876045253317
    "E9" = 2.7723311546840961
    "F9" = 2.7707850281550179
    "G9" = 2.7781972628598397
    "H9" = 2.7786318216756341
    "B10" = 3.7541302235179783
    "C10" = 3.7769532955041467
    "D10" = 3.7825873093949824
    "E10" = 3.7941176470588238
    "F10" = 3.7827095064590921
    "G10" = 3.7796130250117979
    "H10" = 3.7791442480143482
}
foreach ($cellRef in $numericValues.Keys) {
    $ws.Range($cellRef).Value = $numericValues[$cellRef]
}

# 2. Column A (rows 2-5 and 7-9) switches from numeric bucket labels to text
#    labels (e.g. 0.25 -> "0.25", 1 -> "1.0", 2 -> "2.0", 3 -> "3.0") while
#    keeping each cell's existing style (border/alignment) untouched. Writing
#    through a helper cell + PasteSpecial(values) avoids Excel's normal
#    "looks like a number" auto-coercion reformatting the cell style.
$textValues = [ordered]@{
    "A2" = "0.25"
    "A3" = "0.5"
    "A4" = "1.0"
    "A5" = "1.5"
    "A7" = "2.0"
    "A8" = "2.5"
    "A9" = "3.0"
}
foreach ($cellRef in $textValues.Keys) {
    $ws.Range("Z1").Value = "'" + $textValues[$cellRef]
    $ws.Range("Z1").Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
}
$ws.Range("Z1").Clear()

# 3. Match the saved selection state from the authored workbook.
$ws.Range("I1").Select()
